$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 8382.833000000001
$ws.Range("I64").Value = 7832.6665
$ws.Range("K64").Value = 7832.6665
$ws.Range("M64").Value = -7584.6665
$ws.Range("H67").Value = 8382.833000000001
$ws.Range("I67").Value = 7832.6665
$ws.Range("K67").Value = 7832.6665
$ws.Range("M67").Value = -6974.6665

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 66.57143000000001
$ws.Range("I5").Value = 53.2
$ws.Range("K5").Value = 53.2
$ws.Range("M5").Value = 58.8
$ws.Range("H38").Value = 21640.334
$ws.Range("I38").Value = 950
$ws.Range("K38").Value = 950
$ws.Range("M38").Value = -483
$ws.Range("H45").Value = 2310
$ws.Range("I45").Value = 1598.6666
$ws.Range("K45").Value = 1598.6666
$ws.Range("M45").Value = -1221.6666
$ws.Range("H61").Value = 5331.6665
$ws.Range("I61").Value = 3747.1
$ws.Range("K61").Value = 3747.1
$ws.Range("M61").Value = -3535.1
$ws.Range("H63").Value = 8860.556
$ws.Range("I63").Value = 1449
$ws.Range("J63").Value = 18125
$ws.Range("K63").Value = 1449
$ws.Range("L63").Value = 18125
$ws.Range("M63").Value = -763
$ws.Range("N63").Value = -19497
$ws.Range("H66").Value = 8860.556
$ws.Range("I66").Value = 1449
$ws.Range("J66").Value = 18125
$ws.Range("K66").Value = 7245
$ws.Range("L66").Value = 90625
$ws.Range("M66").Value = -3813
$ws.Range("N66").Value = -97489
$ws.Range("H122").Value = 999.8333
$ws.Range("I122").Value = 999.8333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2999.4999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -549.4998999999998
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 5331.6665
$ws.Range("I136").Value = 3747.1
$ws.Range("K136").Value = 11241.3
$ws.Range("M136").Value = -8691.299999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 66.57143000000001
$ws.Range("I4").Value = 53.2
$ws.Range("K4").Value = 53.2
$ws.Range("M4").Value = 61.8
$ws.Range("H35").Value = 99999
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H82").Value = 32089.125
$ws.Range("I82").Value = 7800.75
$ws.Range("J82").Value = 56377.5
$ws.Range("K82").Value = 7800.75
$ws.Range("L82").Value = 56377.5
$ws.Range("M82").Value = -7417.75
$ws.Range("N82").Value = -57143.5
$ws.Range("H85").Value = 32089.125
$ws.Range("I85").Value = 7800.75
$ws.Range("J85").Value = 56377.5
$ws.Range("K85").Value = 7800.75
$ws.Range("L85").Value = 56377.5
$ws.Range("M85").Value = -6474.75
$ws.Range("N85").Value = -59029.5
$ws.Range("H86").Value = 5681.6665
$ws.Range("I86").Value = 4818
$ws.Range("K86").Value = 4818
$ws.Range("M86").Value = -3695
$ws.Range("H89").Value = 5681.6665
$ws.Range("I89").Value = 4818
$ws.Range("K89").Value = 24090
$ws.Range("M89").Value = -18474
$ws.Range("H99").Value = 2600
$ws.Range("I99").Value = 2600
$ws.Range("K99").Value = 2600
$ws.Range("M99").Value = -1102
$ws.Range("H103").Value = 8839.125
$ws.Range("I103").Value = 4499
$ws.Range("J103").Value = 9459.143
$ws.Range("K103").Value = 4499
$ws.Range("L103").Value = 9459.143
$ws.Range("M103").Value = -3327
$ws.Range("N103").Value = -11803.143

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 465.41177
$ws.Range("I7").Value = 88.5
$ws.Range("K7").Value = 88.5
$ws.Range("M7").Value = 24.5
$ws.Range("H22").Value = 199.75
$ws.Range("I22").Value = 199.75
$ws.Range("K22").Value = 199.75
$ws.Range("M22").Value = 150.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 663
$ws.Range("I112").Value = 640.5
$ws.Range("K112").Value = 1921.5
$ws.Range("M112").Value = -813.5
$ws.Range("H117").Value = 157
$ws.Range("J117").Value = 163.09091
$ws.Range("L117").Value = 489.27273
$ws.Range("N117").Value = -7373.27273

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 964.3333
$ws.Range("I31").Value = 964.3333
$ws.Range("K31").Value = 964.3333
$ws.Range("M31").Value = -672.3333
$ws.Range("H37").Value = 964.3333
$ws.Range("I37").Value = 964.3333
$ws.Range("K37").Value = 964.3333
$ws.Range("M37").Value = -687.3333
$ws.Range("H80").Value = 3500
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 3500
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 1678.3636
$ws.Range("J32").Value = 2799.6667
$ws.Range("L32").Value = 2799.6667
$ws.Range("N32").Value = -3433.6667
$ws.Range("H46").Value = 5998
$ws.Range("I46").Value = 5000
$ws.Range("J46").Value = 6663.3335
$ws.Range("K46").Value = 5000
$ws.Range("L46").Value = 6663.3335
$ws.Range("M46").Value = -4812
$ws.Range("N46").Value = -7039.3335
$ws.Range("H68").Value = 5250
$ws.Range("I68").Value = 3666.6667
$ws.Range("K68").Value = 3666.6667
$ws.Range("M68").Value = -2917.6667
$ws.Range("H71").Value = 5250
$ws.Range("I71").Value = 3666.6667
$ws.Range("K71").Value = 18333.3335
$ws.Range("M71").Value = -14589.3335
$ws.Range("H82").Value = 3918.4546
$ws.Range("I82").Value = 625
$ws.Range("J82").Value = 4650.3335
$ws.Range("K82").Value = 625
$ws.Range("L82").Value = 4650.3335
$ws.Range("M82").Value = -264
$ws.Range("N82").Value = -5372.3335
$ws.Range("H85").Value = 3918.4546
$ws.Range("I85").Value = 625
$ws.Range("J85").Value = 4650.3335
$ws.Range("K85").Value = 625
$ws.Range("L85").Value = 4650.3335
$ws.Range("M85").Value = 623
$ws.Range("N85").Value = -7146.3335
$ws.Range("H122").Value = 4918.8
$ws.Range("I122").Value = 4918.8
$ws.Range("K122").Value = 14756.4
$ws.Range("M122").Value = -12306.4
$ws.Range("H132").Value = 2464.1667
$ws.Range("I132").Value = 2464.1667
$ws.Range("K132").Value = 7392.500100000001
$ws.Range("M132").Value = -4862.500100000001
$ws.Range("H136").Value = 1997.6
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10375
$ws.Range("H65").Value = 10375
$ws.Range("H100").Value = 630.625
$ws.Range("I100").Value = 630.625
$ws.Range("K100").Value = 1261.25
$ws.Range("M100").Value = -720.25
